# Insert a new weekly price-report row just before the existing row 76
# (Terminal Hortofrutícola Agro Chillán - Repollo), pushing the old rows
# 76-147 down to 77-148. This mirrors the author's "Fruta / hortaliza,
# semanal" commit, which adds one more weekly observation to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 76 (and everything below it) down by one row.
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly record.
$ws.Cells.Item(76, 1).Value  = 7
$ws.Cells.Item(76, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(76, 3).Value  = 'Ñuble'
$ws.Cells.Item(76, 4).Value  = 44512
$ws.Cells.Item(76, 5).Value  = 16
$ws.Cells.Item(76, 6).Value  = 100112006
$ws.Cells.Item(76, 7).Value  = 'Repollo'
$ws.Cells.Item(76, 8).Value  = 'Crespo record'
$ws.Cells.Item(76, 9).Value  = 'Primera'
$ws.Cells.Item(76, 10).Value = 300
$ws.Cells.Item(76, 11).Value = 600
$ws.Cells.Item(76, 12).Value = 700
$ws.Cells.Item(76, 13).Value = 650
$ws.Cells.Item(76, 14).Value = '$/unidad'
$ws.Cells.Item(76, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(76, 16).Value = 650
$ws.Cells.Item(76, 17).Value = 1
$ws.Cells.Item(76, 18).Value = 'Hortaliza'

# Make sure the date cell keeps the original date number format used
# throughout column D.
$ws.Cells.Item(76, 4).NumberFormat = $ws.Cells.Item(77, 4).NumberFormat
